$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 97; existing rows 97-199 shift down to 99-201.
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(97).Insert()

# Populate the first new row (97) - Primera quality entry for 2021-10-15 (serial 44484)
$ws.Cells.Item(97,1).Value = 4
$ws.Cells.Item(97,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97,3).Value = "Los Lagos"
$ws.Cells.Item(97,4).Value = 44484
$ws.Cells.Item(97,5).Value = 10
$ws.Cells.Item(97,6).Value = 100112008
$ws.Cells.Item(97,7).Value = "Coliflor"
$ws.Cells.Item(97,8).Value = "Sin especificar"
$ws.Cells.Item(97,9).Value = "Primera"
$ws.Cells.Item(97,10).Value = 700
$ws.Cells.Item(97,11).Value = 1200
$ws.Cells.Item(97,12).Value = 1200
$ws.Cells.Item(97,13).Value = 1200
$ws.Cells.Item(97,14).Value = "`$/unidad"
$ws.Cells.Item(97,15).Value = "Región Metropolitana"
$ws.Cells.Item(97,16).Value = 1200
$ws.Cells.Item(97,17).Value = 1
$ws.Cells.Item(97,18).Value = "Hortaliza"

# Populate the second new row (98) - Segunda quality entry for the same date
$ws.Cells.Item(98,1).Value = 4
$ws.Cells.Item(98,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98,3).Value = "Los Lagos"
$ws.Cells.Item(98,4).Value = 44484
$ws.Cells.Item(98,5).Value = 10
$ws.Cells.Item(98,6).Value = 100112008
$ws.Cells.Item(98,7).Value = "Coliflor"
$ws.Cells.Item(98,8).Value = "Sin especificar"
$ws.Cells.Item(98,9).Value = "Segunda"
$ws.Cells.Item(98,10).Value = 700
$ws.Cells.Item(98,11).Value = 1100
$ws.Cells.Item(98,12).Value = 1100
$ws.Cells.Item(98,13).Value = 1100
$ws.Cells.Item(98,14).Value = "`$/unidad"
$ws.Cells.Item(98,15).Value = "Región Metropolitana"
$ws.Cells.Item(98,16).Value = 1100
$ws.Cells.Item(98,17).Value = 1
$ws.Cells.Item(98,18).Value = "Hortaliza"
